function Set-TextCell($ws, $addr, $val) {
    $c = $ws.Range($addr)
    if ($val -match '^[+-]?\d+(\.\d+)?$') {
        $c.Value = "'" + $val
        $c.Style = "Normal"
    } else {
        $c.Value = $val
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextCell $ws "D2" "60.928.82"
Set-TextCell $ws "E2" "  -0.15%  "
Set-TextCell $ws "D3" "2.888.12"
Set-TextCell $ws "E3" "  -1.26%  "
Set-TextCell $ws "D4" "1.00"
Set-TextCell $ws "E4" "  +0.03%  "
Set-TextCell $ws "D5" "587.96"
Set-TextCell $ws "E5" "  -0.52%  "
Set-TextCell $ws "D6" "138.44"
Set-TextCell $ws "E6" "  -5.85%  "
Set-TextCell $ws "E7" "  -0.03%  "
Set-TextCell $ws "E8" "  -3.05%  "
Set-TextCell $ws "D9" "6.97"
Set-TextCell $ws "E9" "  +0.46%  "
Set-TextCell $ws "E10" "  -4.73%  "
Set-TextCell $ws "D11" "0.426"
Set-TextCell $ws "E11" "  -3.24%  "
Set-TextCell $ws "D12" "0.0000217"
Set-TextCell $ws "E12" "  -4.10%  "
Set-TextCell $ws "D13" "32.18"
Set-TextCell $ws "E13" "  -4.54%  "
Set-TextCell $ws "E14" "  -0.48%  "
Set-TextCell $ws "D15" "3.367.16"
Set-TextCell $ws "E15" "  -1.27%  "
Set-TextCell $ws "D16" "60.843.47"
Set-TextCell $ws "E16" "  -0.20%  "
Set-TextCell $ws "D17" "2.893.16"
Set-TextCell $ws "E17" "  -1.10%  "
Set-TextCell $ws "D18" "6.47"
Set-TextCell $ws "E18" "  -3.57%  "
Set-TextCell $ws "D19" "424.15"
Set-TextCell $ws "E19" "  -1.74%  "
Set-TextCell $ws "D20" "13.17"
Set-TextCell $ws "E20" "  -1.73%  "
Set-TextCell $ws "D21" "0.652"
Set-TextCell $ws "E21" "  -4.04%  "
Set-TextCell $ws "D22" "6.91"
Set-TextCell $ws "E22" "  -2.51%  "
Set-TextCell $ws "D23" "79.80"
Set-TextCell $ws "E23" "  -1.96%  "
Set-TextCell $ws "D24" "10.32"
Set-TextCell $ws "E24" "  -5.63%  "
Set-TextCell $ws "E25" "  -0.03%  "
Set-TextCell $ws "E26" "  -7.35%  "
Set-TextCell $ws "D27" "11.35"
Set-TextCell $ws "E27" "  -4.68%  "
Set-TextCell $ws "E28" "  -3.45%  "
Set-TextCell $ws "D29" "2.06"
Set-TextCell $ws "E29" "  -9.22%  "
Set-TextCell $ws "D30" "6.60"
Set-TextCell $ws "E30" "  -5.98%  "
Set-TextCell $ws "E31" "  +0.02%  "
Set-TextCell $ws "D32" "25.55"
Set-TextCell $ws "E32" "  -4.23%  "
Set-TextCell $ws "D34" "0.0₃0836"
Set-TextCell $ws "E34" "  -3.47%  "
Set-TextCell $ws "E35" "  -4.41%  "
Set-TextCell $ws "D36" "5.42"
Set-TextCell $ws "E36" "  -4.07%  "
Set-TextCell $ws "D37" "48.95"
Set-TextCell $ws "E37" "  -2.18%  "
Set-TextCell $ws "D38" "2.79"
Set-TextCell $ws "E38" "  -7.40%  "
Set-TextCell $ws "D39" "1.89"
Set-TextCell $ws "E39" "  -4.67%  "
Set-TextCell $ws "D40" "8.31"
Set-TextCell $ws "E40" "  -3.07%  "
Set-TextCell $ws "D41" "0.115"
Set-TextCell $ws "E41" "  -5.52%  "
Set-TextCell $ws "D42" "0.264"
Set-TextCell $ws "E42" "  -6.45%  "
Set-TextCell $ws "D43" "37.75"
Set-TextCell $ws "E43" "  -8.99%  "
Set-TextCell $ws "D44" "2.659.93"
Set-TextCell $ws "E44" "  -1.76%  "
Set-TextCell $ws "D45" "131.11"
Set-TextCell $ws "E45" "  -2.21%  "
Set-TextCell $ws "E46" "  -4.88%  "
Set-TextCell $ws "D47" "350.08"
Set-TextCell $ws "E47" "  -7.49%  "
Set-TextCell $ws "E48" "  +0.03%  "
Set-TextCell $ws "D50" "22.15"
Set-TextCell $ws "E50" "  -7.16%  "
Set-TextCell $ws "E51" "  -4.82%  "
